$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 102, shifting rows 102:210 down to 103:211
$ws.Rows(102).Insert()

# Fill in the new row 102 with the inserted record's data
$ws.Range("A102").Value = 8
$ws.Range("B102").Value = "Terminal La Palmera de La Serena"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 45159
$ws.Range("E102").Value = 4
$ws.Range("F102").Value = 100112052
$ws.Range("G102").Value = "Albahaca"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 800
$ws.Range("K102").Value = 2800
$ws.Range("L102").Value = 3000
$ws.Range("M102").Value = 2900
$ws.Range("N102").Value = "$/paquete"
$ws.Range("O102").Value = "Región de Arica y Parinacota"
$ws.Range("P102").Value = 2900
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = "Hortaliza"
